$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.724.90"
Set-TextValue $ws.Range("E2") "  -0.64%  "

Set-TextValue $ws.Range("D3") "3.797.85"
Set-TextValue $ws.Range("E3") "  -1.50%  "

Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.02%  "

Set-TextValue $ws.Range("D5") "703.92"
Set-TextValue $ws.Range("E5") "  +1.08%  "

Set-TextValue $ws.Range("D6") "170.11"
Set-TextValue $ws.Range("E6") "  -1.79%  "

Set-TextValue $ws.Range("D7") "3.797.09"
Set-TextValue $ws.Range("E7") "  -1.49%  "

Set-TextValue $ws.Range("E8") "  +0.04%  "

Set-TextValue $ws.Range("E9") "  -1.20%  "

Set-TextValue $ws.Range("D10") "0.160"
Set-TextValue $ws.Range("E10") "  -2.07%  "

Set-TextValue $ws.Range("D11") "7.36"
Set-TextValue $ws.Range("E11") "  +0.26%  "

Set-TextValue $ws.Range("E12") "  -1.49%  "

Set-TextValue $ws.Range("E13") "  -2.17%  "

Set-TextValue $ws.Range("D14") "36.04"
Set-TextValue $ws.Range("E14") "  -1.56%  "

Set-TextValue $ws.Range("D15") "4.442.18"
Set-TextValue $ws.Range("E15") "  -1.42%  "

Set-TextValue $ws.Range("D16") "3.794.26"
Set-TextValue $ws.Range("E16") "  -1.82%  "

Set-TextValue $ws.Range("D17") "70.695.76"
Set-TextValue $ws.Range("E17") "  -0.68%  "

Set-TextValue $ws.Range("E18") "  +0.05%  "

Set-TextValue $ws.Range("D19") "7.13"
Set-TextValue $ws.Range("E19") "  -1.81%  "

Set-TextValue $ws.Range("E20") "  -2.28%  "

Set-TextValue $ws.Range("D21") "497.75"
Set-TextValue $ws.Range("E21") "  +0.28%  "

Set-TextValue $ws.Range("E22") "  -4.90%  "

Set-TextValue $ws.Range("D23") "0.726"
Set-TextValue $ws.Range("E23") "  +0.21%  "

Set-TextValue $ws.Range("D24") "84.80"
Set-TextValue $ws.Range("E24") "  -0.27%  "

Set-TextValue $ws.Range("E25") "  -1.51%  "

Set-TextValue $ws.Range("D26") "12.09"
Set-TextValue $ws.Range("E26") "  -2.22%  "

Set-TextValue $ws.Range("D27") "10.41"
Set-TextValue $ws.Range("E27") "  -1.83%  "

Set-TextValue $ws.Range("D28") "3.949.29"
Set-TextValue $ws.Range("E28") "  -1.50%  "

Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  -0.06%  "

Set-TextValue $ws.Range("D30") "2.04"
Set-TextValue $ws.Range("E30") "  -5.10%  "

Set-TextValue $ws.Range("D31") "3.07"
Set-TextValue $ws.Range("E31") "  -2.84%  "

Set-TextValue $ws.Range("E32") "  -4.27%  "

Set-TextValue $ws.Range("E33") "  -4.26%  "

Set-TextValue $ws.Range("D34") "29.02"
Set-TextValue $ws.Range("E34") "  -2.52%  "

Set-TextValue $ws.Range("D35") "0.173"
Set-TextValue $ws.Range("E35") "  -2.90%  "

Set-TextValue $ws.Range("E36") "  +0.03%  "

Set-TextValue $ws.Range("D37") "3.769.13"

Set-TextValue $ws.Range("D38") "9.07"
Set-TextValue $ws.Range("E38") "  -2.60%  "

Set-TextValue $ws.Range("E39") "  -3.85%  "

Set-TextValue $ws.Range("E40") "  +1.47%  "

Set-TextValue $ws.Range("D41") "2.31"
Set-TextValue $ws.Range("E41") "  -3.04%  "

Set-TextValue $ws.Range("D42") "5.93"
Set-TextValue $ws.Range("E42") "  -2.26%  "

Set-TextValue $ws.Range("D43") "3.27"
Set-TextValue $ws.Range("E43") "  -4.14%  "

Set-TextValue $ws.Range("E44") "  -0.04%  "

Set-TextValue $ws.Range("E45") "  +0.25%  "

Set-TextValue $ws.Range("D46") "0.000322"
Set-TextValue $ws.Range("E46") "  +4.94%  "

Set-TextValue $ws.Range("D47") "164.72"
Set-TextValue $ws.Range("E47") "  +0.19%  "

Set-TextValue $ws.Range("D48") "425.26"
Set-TextValue $ws.Range("E48") "  +1.19%  "

Set-TextValue $ws.Range("E49") "  +0.12%  "

Set-TextValue $ws.Range("D50") "8.59"
Set-TextValue $ws.Range("E50") "  -1.10%  "

Set-TextValue $ws.Range("E51") "  -1.98%  "
